$wb = $excel.ActiveWorkbook

# --- Sheet "Data": replace the 5 identical rows with a single header-like row of new data ---
$dataSheet = $wb.Worksheets.Item("Data")

# Clear existing contents (rows 1-5, cols A-C) then write new single row of values
$dataSheet.Range("A1:C5").ClearContents()

$dataSheet.Range("A1").Value = "022A-9661"
$dataSheet.Range("B1").Value = "Jiro Miko Vinas"
$dataSheet.Range("C1").Value = "jiromikovinas@gmail.com"
$dataSheet.Range("D1").Value = "Male"
$dataSheet.Range("E1").Value = "BSIT"
# Format as text first so the leading zero in the mobile number is preserved
$dataSheet.Range("F1").NumberFormat = "@"
$dataSheet.Range("F1").Value = "09518583657"
$dataSheet.Range("G1").Value = "Lucena City"

# --- Sheet "Old Students": append a new row 4 ---
$oldStudents = $wb.Worksheets.Item("Old Students")

$oldStudents.Range("A4").Value = "022A-9661"
$oldStudents.Range("B4").Value = "jiromiko"
$oldStudents.Range("C4").Value = "BSIT-1A"
